$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the "inspired by" + " the first or first few mansions of "
#    runs into a single run (same combined text, empty run formatting).
# ---------------------------------------------------------------------
$tailRng = $d.Content.Duplicate
$tailRng.Find.Execute(" the first or first few mansions of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailRng.Delete()

$headRng = $d.Content.Duplicate
$headRng.Find.Execute("inspired by", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headRng.Collapse(0)
$headRng.InsertAfter(" the first or first few mansions of ")

# nudge formatting on/off to coax the run back to an explicit, empty
# <w:rPr/> on the freshly-merged run (matches the surrounding runs).
$mergedRng = $d.Content.Duplicate
$mergedRng.Find.Execute("inspired by the first or first few mansions of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergedRng.Font.Bold = $true
$mergedRng.Font.Bold = $false

# ---------------------------------------------------------------------
# 2. Split "area/enemy/item variety, difficulty curve" so a new
#    "boss fight(s)" run follows it (comma + space stay on the first
#    run).
# ---------------------------------------------------------------------
$curveRng = $d.Content.Duplicate
$curveRng.Find.Execute("area/enemy/item variety, difficulty curve", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$curveRng.Collapse(0)
$curveRng.InsertAfter(", boss fight(s)")

$bossRng = $d.Content.Duplicate
$bossRng.Find.Execute("boss fight(s)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bossRng.Font.Bold = $true
$bossRng.Font.Bold = $false

# ---------------------------------------------------------------------
# 3. Flip "overflowPunct" (Word's hanging-punctuation toggle) on the
#    Normal style's paragraph formatting from false to true.
# ---------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $true
